# "basic project info stuff"
#
# 1. Agenda slide ("Agenda", slide 2): the second content placeholder
#    ("Text Placeholder 3") had two empty bulleted paragraphs. Put
#    "Photos" into the second paragraph and leave a fresh empty bulleted
#    paragraph after it (so the placeholder ends up with three
#    paragraphs: blank, "Photos", blank).
# 2. Insert two new Title-Slide-layout slides ("Introductions" and
#    "Photos") right before the closing "THANK YOU" slide.

$p = $ppt.ActivePresentation

$agendaSlide = $p.Slides.Item(2)
$photosShape = $agendaSlide.Shapes.Item(3)
$tr = $photosShape.TextFrame.TextRange
$photosPara = $tr.Paragraphs(2, 1)
$photosPara.Text = "Photos"
$photosPara = $tr.Paragraphs(2, 1)
[void]$photosPara.InsertAfter("`r")

$introSlide = $p.Slides.Add(3, 1)
$introSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Introductions"

$photosSlide = $p.Slides.Add(4, 1)
$photosSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Photos"
